$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.559.50"
$ws.Range("E2").Value = "  -1.27%  "
$ws.Range("D3").Value = "1.670.93"
$ws.Range("E3").Value = "  -2.25%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3953"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3945"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.67%  "
$ws.Range("E9").Value = "  -0.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.396"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "50.74"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08642"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.48"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.317"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001318"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.693"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.96%  "
$ws.Range("D17").Value = "1.696.30"
$ws.Range("E17").Value = "  +1.88%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.95"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07011"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "21.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.073"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.98%  "
$ws.Range("E22").Value = "  -0.39%  "
$ws.Range("E23").Value = "  -4.00%  "
$ws.Range("D24").Value = "24.571.19"
$ws.Range("E24").Value = "  -1.22%  "
$ws.Range("E25").Value = "  +0.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.762"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.835"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -8.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "160.05"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "145.83"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.294"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.520"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +10.67%  "
$ws.Range("D33").Value = "1.861.40"
$ws.Range("E33").Value = "  +0.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.03081"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08255"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.81%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.903"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2807"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9936"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.09632"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.517"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.29"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7888"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.56"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.69%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.54"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.86%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.563"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.91%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.7095"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.176"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.40%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.08631"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.002"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.324"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "138.28"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.75%  "
